$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.365258812904358
$ws.Range("B1").Value = 2.796512126922607
$ws.Range("C1").Value = 3.663173913955688
$ws.Range("D1").Value = 3.588552951812744
$ws.Range("E1").Value = 1.15246307849884
